$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet (tab) name to reflect new "through" date
$ws.Name = "Through 2022-06-29"

# Update the label for the June row (shared string used in A7)
$ws.Range("A7").Value = "June (through 06-29)"

# Update June row (row 7) values for columns C..I (B unchanged)
$ws.Range("C7").Value = 39
$ws.Range("D7").Value = 71
$ws.Range("E7").Value = 57
$ws.Range("F7").Value = 46
$ws.Range("G7").Value = 112
$ws.Range("H7").Value = 123
$ws.Range("I7").Value = 139

# Update Total row (row 8) values for columns C..I (B unchanged)
$ws.Range("C8").Value = 248
$ws.Range("D8").Value = 387
$ws.Range("E8").Value = 352
$ws.Range("F8").Value = 250
$ws.Range("G8").Value = 470
$ws.Range("H8").Value = 754
$ws.Range("I8").Value = 802

$wb.Save()
